$d = $word.ActiveDocument

$replacements = @(
    @{old="909×2="; new="846×3="},
    @{old="885×3="; new="901×6="},
    @{old="791×8="; new="929×8="},
    @{old="621×8="; new="936×4="},
    @{old="964×3="; new="304×5="},
    @{old="118×7="; new="922×8="},
    @{old="730×8="; new="182×6="},
    @{old="776×5="; new="990×8="},
    @{old="436×4="; new="272×2="},
    @{old="934×9="; new="251×4="},
    @{old="940×7="; new="601×4="},
    @{old="911×4="; new="497×7="},
    @{old="338×8="; new="362×3="},
    @{old="678×7="; new="573×8="},
    @{old="415×6="; new="428×9="},
    @{old="856×2="; new="980×6="},
    @{old="830×7="; new="941×3="},
    @{old="334×6="; new="231×8="},
    @{old="142×4="; new="485×8="},
    @{old="541×3="; new="935×8="},
    @{old="516×3="; new="269×7="},
    @{old="562×2="; new="616×8="},
    @{old="805×7="; new="482×6="},
    @{old="668×8="; new="971×2="},
    @{old="281×8="; new="906×5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
